# Refresh the "cryptos" price list (scheduled GitHub Actions data pull).
# Most rows only get new Price (D) / Volume(1h) (E) figures; rows 41-42
# additionally swap Coin (B) and Link (C) because Maker and VeChain traded
# ranking positions between the two snapshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ForceText: some "Price" strings look like plain numbers (e.g. "228.46").
# Excel's Range.Value setter auto-converts such strings to real numbers,
# which would lose formatting (e.g. "0.140" -> 0.14) and change the cell's
# stored type. Flip the cell to the "@" (Text) format for the assignment,
# then restore its original format, to keep these as text - same as the
# source data.
$updates = @(
    @{ Addr = "D2"; Value = "38.791.81"; ForceText = $false },
    @{ Addr = "E2"; Value = "  +2.64%  "; ForceText = $false },
    @{ Addr = "D3"; Value = "2.089.26"; ForceText = $false },
    @{ Addr = "E3"; Value = "  +1.98%  "; ForceText = $false },
    @{ Addr = "E4"; Value = "  +0.06%  "; ForceText = $false },
    @{ Addr = "D5"; Value = "228.46"; ForceText = $true },
    @{ Addr = "E5"; Value = "  +0.24%  "; ForceText = $false },
    @{ Addr = "E6"; Value = "  +0.75%  "; ForceText = $false },
    @{ Addr = "D7"; Value = "60.24"; ForceText = $true },
    @{ Addr = "E7"; Value = "  +0.66%  "; ForceText = $false },
    @{ Addr = "E8"; Value = "  +0.02%  "; ForceText = $false },
    @{ Addr = "D9"; Value = "0.384"; ForceText = $true },
    @{ Addr = "E9"; Value = "  +1.75%  "; ForceText = $false },
    @{ Addr = "D10"; Value = "0.0841"; ForceText = $true },
    @{ Addr = "E10"; Value = "  +0.79%  "; ForceText = $false },
    @{ Addr = "E11"; Value = "  -0.24%  "; ForceText = $false },
    @{ Addr = "D12"; Value = "2.402.15"; ForceText = $false },
    @{ Addr = "E12"; Value = "  +2.14%  "; ForceText = $false },
    @{ Addr = "D13"; Value = "14.99"; ForceText = $true },
    @{ Addr = "E13"; Value = "  +3.98%  "; ForceText = $false },
    @{ Addr = "D14"; Value = "21.89"; ForceText = $true },
    @{ Addr = "E14"; Value = "  +1.88%  "; ForceText = $false },
    @{ Addr = "D15"; Value = "0.797"; ForceText = $true },
    @{ Addr = "E15"; Value = "  +4.25%  "; ForceText = $false },
    @{ Addr = "E16"; Value = "  -0.72%  "; ForceText = $false },
    @{ Addr = "D17"; Value = "2.097.43"; ForceText = $false },
    @{ Addr = "E17"; Value = "  +2.02%  "; ForceText = $false },
    @{ Addr = "D18"; Value = "38.744.46"; ForceText = $false },
    @{ Addr = "E18"; Value = "  +2.51%  "; ForceText = $false },
    @{ Addr = "D19"; Value = "71.55"; ForceText = $true },
    @{ Addr = "E19"; Value = "  +2.90%  "; ForceText = $false },
    @{ Addr = "E20"; Value = "  +2.02%  "; ForceText = $false },
    @{ Addr = "D21"; Value = "0.0₃0838"; ForceText = $false },
    @{ Addr = "E21"; Value = "  +0.99%  "; ForceText = $false },
    @{ Addr = "D22"; Value = "227.22"; ForceText = $true },
    @{ Addr = "E22"; Value = "  +2.11%  "; ForceText = $false },
    @{ Addr = "E23"; Value = "  -0.41%  "; ForceText = $false },
    @{ Addr = "E24"; Value = "  -0.36%  "; ForceText = $false },
    @{ Addr = "D25"; Value = "2.34"; ForceText = $true },
    @{ Addr = "E25"; Value = "  +2.46%  "; ForceText = $false },
    @{ Addr = "D26"; Value = "170.78"; ForceText = $true },
    @{ Addr = "E26"; Value = "  +1.13%  "; ForceText = $false },
    @{ Addr = "D27"; Value = "9.51"; ForceText = $true },
    @{ Addr = "E27"; Value = "  +1.96%  "; ForceText = $false },
    @{ Addr = "D28"; Value = "0.140"; ForceText = $true },
    @{ Addr = "E28"; Value = "  +9.66%  "; ForceText = $false },
    @{ Addr = "D29"; Value = "1.47"; ForceText = $true },
    @{ Addr = "E29"; Value = "  +13.07%  "; ForceText = $false },
    @{ Addr = "D30"; Value = "19.16"; ForceText = $true },
    @{ Addr = "E30"; Value = "  +1.96%  "; ForceText = $false },
    @{ Addr = "D31"; Value = "0.120"; ForceText = $true },
    @{ Addr = "E31"; Value = "  +0.77%  "; ForceText = $false },
    @{ Addr = "E32"; Value = "  +5.73%  "; ForceText = $false },
    @{ Addr = "E33"; Value = "  +2.68%  "; ForceText = $false },
    @{ Addr = "D34"; Value = "4.69"; ForceText = $true },
    @{ Addr = "E34"; Value = "  +3.59%  "; ForceText = $false },
    @{ Addr = "D35"; Value = "0.0609"; ForceText = $true },
    @{ Addr = "E35"; Value = "  +1.07%  "; ForceText = $false },
    @{ Addr = "D36"; Value = "6.48"; ForceText = $true },
    @{ Addr = "E36"; Value = "  -0.10%  "; ForceText = $false },
    @{ Addr = "E37"; Value = "  +1.24%  "; ForceText = $false },
    @{ Addr = "D38"; Value = "3.57"; ForceText = $true },
    @{ Addr = "E38"; Value = "  +2.36%  "; ForceText = $false },
    @{ Addr = "E39"; Value = "  +0.18%  "; ForceText = $false },
    @{ Addr = "D40"; Value = "18.09"; ForceText = $true },
    @{ Addr = "E40"; Value = "  -1.93%  "; ForceText = $false },
    @{ Addr = "B41"; Value = "VeChain"; ForceText = $false },
    @{ Addr = "C41"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; ForceText = $false },
    @{ Addr = "D41"; Value = "0.0226"; ForceText = $true },
    @{ Addr = "E41"; Value = "  +4.75%  "; ForceText = $false },
    @{ Addr = "B42"; Value = "Maker"; ForceText = $false },
    @{ Addr = "C42"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; ForceText = $false },
    @{ Addr = "D42"; Value = "1.542.19"; ForceText = $false },
    @{ Addr = "E42"; Value = "  +1.10%  "; ForceText = $false },
    @{ Addr = "D43"; Value = "100.72"; ForceText = $true },
    @{ Addr = "E43"; Value = "  +3.19%  "; ForceText = $false },
    @{ Addr = "E44"; Value = "  -0.73%  "; ForceText = $false },
    @{ Addr = "D45"; Value = "0.0924"; ForceText = $true },
    @{ Addr = "E45"; Value = "  +3.73%  "; ForceText = $false },
    @{ Addr = "E46"; Value = "  +8.41%  "; ForceText = $false },
    @{ Addr = "E47"; Value = "  +1.57%  "; ForceText = $false },
    @{ Addr = "D48"; Value = "4.11"; ForceText = $true },
    @{ Addr = "E48"; Value = "  -2.22%  "; ForceText = $false },
    @{ Addr = "E49"; Value = "  +2.48%  "; ForceText = $false },
    @{ Addr = "D50"; Value = "2.96"; ForceText = $true },
    @{ Addr = "E50"; Value = "  +0.66%  "; ForceText = $false },
    @{ Addr = "D51"; Value = "2.290.29"; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Addr)
    if ($u.ForceText) {
        $savedFormat = $range.NumberFormat
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.NumberFormat = $savedFormat
    } else {
        $range.Value = $u.Value
    }
}
